# Registrar Usuario - fin con fallos
# Adds the new "Ana Perez" user row (row 5) to the Usuarios sheet, finishing
# the "add user" flow (still pending: fixing bugs related to that new
# user's subsequent login).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Plain text columns (names / sex / observation) - safe to assign directly,
# Excel won't reinterpret these as numbers.
$ws.Range("B5").Value = "Ana"
$ws.Range("C5").Value = "Perez"
$ws.Range("E5").Value = "M"
$ws.Range("K5").Value = "mantener"

# Numeric-looking values that must be stored as TEXT (id, password, edad,
# altura, peso) - build each as a formula returning text, then paste back
# as a value so the cell keeps General formatting/style (no NumberFormat
# change) while the stored type stays a string, matching how this data was
# entered on the source sheet.
$textCells = @(
    @{ Addr = "A5"; Num = 3 },
    @{ Addr = "D5"; Num = 3 },
    @{ Addr = "F5"; Num = 21 },
    @{ Addr = "G5"; Num = 168 },
    @{ Addr = "H5"; Num = 55 }
)
foreach ($item in $textCells) {
    $cell = $ws.Range($item.Addr)
    $cell.Formula = "=TEXT(" + $item.Num + ",""0"")"
    $cell.Copy()
    $cell.PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteValues)
}

# True numeric columns (actividad / patologia) stay as numbers.
$ws.Range("I5").Value = 2
$ws.Range("J5").Value = 0
